$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'30.322.46"
$ws.Range("E2").Value = "'  +0.09%  "

# Row 3
$ws.Range("D3").Value = "'1.931.46"
$ws.Range("E3").Value = "'  +0.19%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'0.7494"
$ws.Range("E5").Value = "'  +4.36%  "

# Row 6
$ws.Range("D6").Value = "'243.54"
$ws.Range("E6").Value = "'  -2.21%  "

# Row 7
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "'  +0.00%  "

# Row 8
$ws.Range("D8").Value = "'0.3186"
$ws.Range("E8").Value = "'  -0.52%  "

# Row 9
$ws.Range("D9").Value = "'27.55"
$ws.Range("E9").Value = "'  -1.10%  "

# Row 10
$ws.Range("D10").Value = "'0.07124"
$ws.Range("E10").Value = "'  +0.61%  "

# Row 11
$ws.Range("D11").Value = "'0.7830"
$ws.Range("E11").Value = "'  -0.74%  "

# Row 12
$ws.Range("D12").Value = "'0.08055"
$ws.Range("E12").Value = "'  +0.89%  "

# Row 13
$ws.Range("D13").Value = "'1.956.84"
$ws.Range("E13").Value = "'  +1.46%  "

# Row 14
$ws.Range("D14").Value = "'5.408"
$ws.Range("E14").Value = "'  +0.60%  "

# Row 15
$ws.Range("D15").Value = "'93.30"
$ws.Range("E15").Value = "'  -1.49%  "

# Row 16
$ws.Range("D16").Value = "'14.62"
$ws.Range("E16").Value = "'  -0.09%  "

# Row 17
$ws.Range("D17").Value = "'30.324.96"
$ws.Range("E17").Value = "'  +0.08%  "

# Row 18
$ws.Range("B18").Value = "'Uniswap"
$ws.Range("C18").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'6.030"
$ws.Range("E18").Value = "'  +4.79%  "

# Row 19
$ws.Range("B19").Value = "'BitcoinCash"
$ws.Range("C19").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'252.51"
$ws.Range("E19").Value = "'  -1.81%  "

# Row 20
$ws.Range("D20").Value = "'0.000007940"
$ws.Range("E20").Value = "'  -1.88%  "

# Row 21
$ws.Range("B21").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "'2.207.86"
$ws.Range("E21").Value = "'  +1.16%  "

# Row 22
$ws.Range("B22").Value = "'Dai"
$ws.Range("C22").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "'  -0.01%  "

# Row 23
$ws.Range("E23").Value = "'  -0.04%  "

# Row 24
$ws.Range("D24").Value = "'6.675"
$ws.Range("E24").Value = "'  -2.12%  "

# Row 25
$ws.Range("D25").Value = "'9.599"
$ws.Range("E25").Value = "'  +0.66%  "

# Row 26
$ws.Range("D26").Value = "'165.46"
$ws.Range("E26").Value = "'  +0.60%  "

# Row 27
$ws.Range("D27").Value = "'19.12"
$ws.Range("E27").Value = "'  +0.13%  "

# Row 28
$ws.Range("D28").Value = "'0.1301"
$ws.Range("E28").Value = "'  +2.48%  "

# Row 29
$ws.Range("D29").Value = "'2.191"
$ws.Range("E29").Value = "'  -3.46%  "

# Row 30
$ws.Range("E30").Value = "'  +0.79%  "

# Row 31
$ws.Range("D31").Value = "'1.564"
$ws.Range("E31").Value = "'  +2.23%  "

# Row 32
$ws.Range("D32").Value = "'4.433"
$ws.Range("E32").Value = "'  +0.69%  "

# Row 33
$ws.Range("D33").Value = "'4.147"
$ws.Range("E33").Value = "'  +0.26%  "

# Row 34
$ws.Range("D34").Value = "'0.05264"
$ws.Range("E34").Value = "'  +2.48%  "

# Row 35
$ws.Range("D35").Value = "'1.321"

# Row 36
$ws.Range("D36").Value = "'0.7599"
$ws.Range("E36").Value = "'  +2.05%  "

# Row 37
$ws.Range("E37").Value = "'  +0.25%  "

# Row 38
$ws.Range("E38").Value = "'  -1.32%  "

# Row 39
$ws.Range("E39").Value = "'  -0.02%  "

# Row 40
$ws.Range("D40").Value = "'6.535"
$ws.Range("E40").Value = "'  +2.52%  "

# Row 41
$ws.Range("D41").Value = "'77.90"
$ws.Range("E41").Value = "'  -0.31%  "

# Row 42
$ws.Range("D42").Value = "'0.4542"
$ws.Range("E42").Value = "'  +0.77%  "

# Row 43
$ws.Range("D43").Value = "'1.976"
$ws.Range("E43").Value = "'  -0.99%  "

# Row 44
$ws.Range("D44").Value = "'0.8433"
$ws.Range("E44").Value = "'  -0.34%  "

# Row 45
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "'  +0.03%  "

# Row 46
$ws.Range("B46").Value = "'EnergySwap"
$ws.Range("C46").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.08"
$ws.Range("E46").Value = "'  +2.82%  "

# Row 47
$ws.Range("B47").Value = "'Aptos"
$ws.Range("C47").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.720"
$ws.Range("E47").Value = "'  +3.81%  "

# Row 48
$ws.Range("D48").Value = "'101.63"
$ws.Range("E48").Value = "'  +0.93%  "

# Row 49
$ws.Range("D49").Value = "'2.105.87"
$ws.Range("E49").Value = "'  +0.86%  "

# Row 50
$ws.Range("D50").Value = "'37.85"
$ws.Range("E50").Value = "'  +2.72%  "

# Row 51
$ws.Range("D51").Value = "'0.1220"
$ws.Range("E51").Value = "'  +6.96%  "
